# Swap the presentation's applied theme from the custom "Integral / Red Violet"
# palette over to the stock "Office" palette (the Design-gallery equivalent of
# clicking the default "Office Theme" tile).
#
# ppt/theme/theme1.xml is the theme actually driving the slide master / slides
# (presentation.xml.rels rId1 -> theme1.xml, slideMaster1.xml.rels rId12 ->
# theme1.xml), so that is the live ColorScheme object PowerPoint exposes via
# SlideMaster.Theme.ThemeColorScheme. Only the 12 scheme colors differ between
# the old and new theme (font scheme + format scheme are identical), so a
# plain ThemeColorScheme.Item(i).RGB walk reproduces the target theme exactly.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# COM RGB values pack as R + G*256 + B*65536, i.e. the reverse byte order of
# the "RRGGBB" hex strings found in <a:srgbClr val="...">.
$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72

# Best-effort: PowerPoint's object model normally treats these as read-only
# (Theme.Name / Design.Name have no settable backing field), but set them
# anyway in case the host persists it.
$theme.Name = "Office Theme"
$colors.Name = "Office"
